$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks like a plain number need NumberFormat forced to Text
# so Excel does not silently convert the string into a numeric value and drop formatting
# such as trailing zeros (e.g. "1.00" -> 1).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.467.46"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "1.602.23"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "212.69"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +7.38%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "26.82"
$ws.Range("E8").Value = "  +8.30%  "
$ws.Range("D9").Value = "43.41"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").Value = "0.0598"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "1.832.30"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").Value = "1.598.63"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "29.510.93"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "0.535"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "63.36"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "242.23"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "3.98"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "2.09"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "154.34"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("D28").Value = "15.28"
$ws.Range("E28").Value = "  +3.34%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.415.80"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "3.09"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "2.79"
$ws.Range("E38").Value = "  +4.75%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").Value = "1.97"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0478"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "52.35"
$ws.Range("E46").Value = "  +20.35%  "
$ws.Range("D47").Value = "65.55"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").Value = "5.29"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "1.743.39"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").Value = "86.37"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  -4.32%  "
